# "feat (excel): store to db"
# Replace the "Dulce" record with the new "Rouge" hire, refresh everyone's
# phone/nid numbers, widen the email column, and leave the cursor on E14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes the new hire (was Dulce, now Rouge) -------------------
$ws.Range("B2").Value = "Rouge"
$ws.Range("C2").Value = "janvierntwali@gmail.com"
$ws.Range("D2").Value = 250724955240
$ws.Range("E2").Value = 1188773399330050
$ws.Range("F2").Value = "designer"
$ws.Range("G2").Value = "15/10/2000"
$ws.Range("H2").Value = "active"

# --- Rows 3-11: same people, refreshed phone/nid numbers -----------------
$ws.Range("B3").Value = "Mara"
$ws.Range("C3").Value = "hashimoto@abc.com"
$ws.Range("D3").Value = 250724955242
$ws.Range("E3").Value = 1188337740373500
$ws.Range("F3").Value = "designer"
$ws.Range("G3").Value = "16/08/2000"
$ws.Range("H3").Value = "active"

$ws.Range("B4").Value = "Philip"
$ws.Range("C4").Value = "Gent@abc.com"
$ws.Range("D4").Value = 250724955241
$ws.Range("E4").Value = 1277422464446730
$ws.Range("F4").Value = "designer"
$ws.Range("G4").Value = "21/05/1994"
$ws.Range("H4").Value = "active"

$ws.Range("B5").Value = "Kathleen"
$ws.Range("C5").Value = "Hanner@abc.com"
$ws.Range("D5").Value = 250724955243
$ws.Range("E5").Value = 1199800819494020
$ws.Range("F5").Value = "developer"
$ws.Range("G5").Value = "15/10/1998"
$ws.Range("H5").Value = "active"

$ws.Range("B6").Value = "Nereida"
$ws.Range("C6").Value = "Magwood@abc.com"
$ws.Range("D6").Value = 250724955244
$ws.Range("E6").Value = 1947300889050020
$ws.Range("F6").Value = "designer"
$ws.Range("G6").Value = "16/08/1999"
$ws.Range("H6").Value = "active"

$ws.Range("B7").Value = "Gaston"
$ws.Range("C7").Value = "Brumm@abc.com"
$ws.Range("D7").Value = 250724955246
$ws.Range("E7").Value = 1483957385758350
$ws.Range("F7").Value = "designer"
$ws.Range("G7").Value = "21/05/2000"
$ws.Range("H7").Value = "active"

$ws.Range("B8").Value = "Etta"
$ws.Range("C8").Value = "Hurn@abc.com"
$ws.Range("D8").Value = 250724955245
$ws.Range("E8").Value = 1947736457285720
$ws.Range("F8").Value = "developer"
$ws.Range("G8").Value = "15/10/1992"
$ws.Range("H8").Value = "inactive"

$ws.Range("B9").Value = "Earlean"
$ws.Range("C9").Value = "Melgar@abc.com"
$ws.Range("D9").Value = 250724955247
$ws.Range("E9").Value = 1893857385003840
$ws.Range("F9").Value = "designer"
$ws.Range("G9").Value = "16/08/1995"
$ws.Range("H9").Value = "inactive"

$ws.Range("B10").Value = "Vincenza"
$ws.Range("C10").Value = "Weiland@abc.com"
$ws.Range("D10").Value = 250724955248
$ws.Range("E10").Value = 7383748728278500
$ws.Range("F10").Value = "developer"
$ws.Range("G10").Value = "21/05/1995"
$ws.Range("H10").Value = "inactive"

$ws.Range("B11").Value = "Fallon"
$ws.Range("C11").Value = "Winward@abc.com"
$ws.Range("D11").Value = 250724955249
$ws.Range("E11").Value = 1199880081551020
$ws.Range("F11").Value = "designer"
$ws.Range("G11").Value = "16/08/1996"
$ws.Range("H11").Value = "inactive"

# --- Column widths: widen "nid" slightly, widen "email" a lot ------------
$ws.Columns.Item(4).ColumnWidth = 15.140625
$ws.Columns.Item(5).ColumnWidth = 35.28515625

# --- Cursor / selection moves to E14 --------------------------------------
[void]$ws.Range("E14").Select()

# --- Page setup: explicit portrait orientation ---------------------------
$ws.PageSetup.Orientation = 1
